$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.707.48"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.902.31"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4974"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3760"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07242"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8932"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07613"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").Value = "1.876.10"
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.439"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008685"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "27.743.86"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.139"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").Value = "2.141.01"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.562"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.841"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.192"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.853"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08906"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.172"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.776"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.224"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7771"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.611"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02067"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.051"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.088"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5496"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05288"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.740"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.443"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1507"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4765"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.621"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06013"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.14%  "
